$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = -21.88
$ws.Range("D6").Value = -8.141999999999999
$ws.Range("A14").Value = -21.644
$ws.Range("B15").Value = 5.45
$ws.Range("A16").Value = -22.121
$ws.Range("D18").Value = -8.485000000000001
$ws.Range("D19").Value = -7.757
$ws.Range("A21").Value = -20.959
$ws.Range("B21").Value = 7.802000000000001
$ws.Range("B22").Value = 7.290999999999999
$ws.Range("A23").Value = -20.856
$ws.Range("B24").Value = 5.524
$ws.Range("A25").Value = -21.589
$ws.Range("A26").Value = -21.012
$ws.Range("B27").Value = 6.44
$ws.Range("B28").Value = 5.583000000000001
$ws.Range("A29").Value = -21.52
$ws.Range("D35").Value = -8.058999999999999
$ws.Range("B36").Value = 6.356
$ws.Range("B39").Value = 6.609
$ws.Range("A40").Value = -20.461
$ws.Range("D44").Value = -7.939
$ws.Range("B45").Value = 5.647
$ws.Range("D47").Value = -7.829000000000001
$ws.Range("B48").Value = 6.238
$ws.Range("B49").Value = 6.470000000000001
$ws.Range("D50").Value = -8.41
$ws.Range("D51").Value = -8.144000000000002
$ws.Range("B52").Value = 5.051
$ws.Range("D52").Value = -8.178000000000001
$ws.Range("A53").Value = -20.723
$ws.Range("B53").Value = 8.217000000000001
$ws.Range("B54").Value = 5.144000000000001
$ws.Range("D55").Value = -8.404999999999999
$ws.Range("A57").Value = -21.527
$ws.Range("B57").Value = 6.327
$ws.Range("D57").Value = -8.17
$ws.Range("D58").Value = -8.430000000000001
$ws.Range("A59").Value = -22.345
$ws.Range("D64").Value = -7.669
$ws.Range("A65").Value = -21.529
$ws.Range("D66").Value = -7.438
$ws.Range("A69").Value = -21.441
$ws.Range("B70").Value = 5.693
$ws.Range("B71").Value = 5.021
$ws.Range("A79").Value = -21.316
$ws.Range("D80").Value = -7.781000000000001
$ws.Range("A83").Value = -22.132
$ws.Range("D83").Value = -8.000999999999999
$ws.Range("B86").Value = 4.937
$ws.Range("B87").Value = 4.528
$ws.Range("B89").Value = 4.903
$ws.Range("A91").Value = -21.033
$ws.Range("D92").Value = -7.283000000000001
$ws.Range("A93").Value = -21.508
$ws.Range("D94").Value = -7.773999999999999
$ws.Range("D96").Value = -7.464
$ws.Range("D97").Value = -7.874
$ws.Range("A100").Value = -22.35
$ws.Range("B101").Value = 5.794
$ws.Range("D101").Value = -7.901999999999999
$ws.Range("A103").Value = -22.055
